$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.192002773284912
$ws.Range("B1").Value = 2.289618253707886
$ws.Range("C1").Value = 6.571616172790527
$ws.Range("D1").Value = 2.313052654266357
$ws.Range("E1").Value = 1.188067674636841
